# Generate Report for Handoff
# Move the localization-status report from "handed back" state to a fresh
# "ready for handoff" state: update the status text, bump the generation
# timestamps, and shrink the now-shorter status columns to fit.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Latest HO Xliff Generate Date / Latest Handback DateTime bump ---
$wsOverview.Range("G2").Value = "2016-08-22 22:58:18"
$wsDeDe.Range("H2").Value     = "2016-08-22 22:58:18"

# --- Latest Handoff Datetime bump (zh-cn) ---
$wsZhCn.Range("H2").Value = "2016-08-22 22:58:13"

# --- Narrow the status columns now that "Ready for handoff" is shorter ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.25   # column E
$wsOverview.Columns.Item(6).ColumnWidth = 16.25   # column F
$wsZhCn.Columns.Item(3).ColumnWidth = 16.25       # column C
$wsDeDe.Columns.Item(3).ColumnWidth = 16.25       # column C
